# Add a new "2022-Q4" quarterly sheet (copied/formatted like the existing
# "2022-Q2" sheet so it keeps the same header/style layout), fill it with the
# new quarter's fund-holdings data, and insert a matching new row into the
# "总计" (totals) summary sheet.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)       # "总计"
$q2    = $wb.Worksheets.Item(2)       # "2022-Q2" (existing template sheet)

# --- 1. Create the new "2022-Q4" sheet right after "总计" -------------------
$q2.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Row 2: 恒生指数基金M类人民币（对冲）份额 — code/name unchanged, stats updated
$q4.Range("D2").Value = "'27.13"
$q4.Range("E2").Value = "'99.24"
$q4.Range("F2").Value = "'8.18"
$q4.Range("G2").Value = "'2.2192"

# Row 3: 建信富时100指数（QDII）人民币A — code/name unchanged, stats updated
$q4.Range("D3").Value = "'0.56"
$q4.Range("E3").Value = "'85.06"
$q4.Range("F3").Value = "'6.04"
$q4.Range("G3").Value = "'0.0338"
$q4.Range("H3").Value = 4

# Row 4: 建信富时100指数（QDII）美元现汇 A
$q4.Range("C4").Value = "建信富时100指数（QDII）美元现汇 A"
$q4.Range("D4").Value = "'0.56"
$q4.Range("E4").Value = "'85.06"
$q4.Range("F4").Value = "'6.04"
$q4.Range("G4").Value = "'0.0338"
$q4.Range("H4").Value = 4

# Row 5: 建信富时100指数（QDII）人民币 C
$q4.Range("C5").Value = "建信富时100指数（QDII）人民币 C"
$q4.Range("D5").Value = "'0.26"
$q4.Range("E5").Value = "'85.06"
$q4.Range("F5").Value = "'6.04"
$q4.Range("G5").Value = "'0.0157"
$q4.Range("H5").Value = 4

# Row 6: 建信富时100指数（QDII）美元现汇 C
$q4.Range("C6").Value = "建信富时100指数（QDII）美元现汇 C"
$q4.Range("D6").Value = "'0.26"
$q4.Range("E6").Value = "'85.06"
$q4.Range("F6").Value = "'6.04"
$q4.Range("G6").Value = "'0.0157"
$q4.Range("H6").Value = 4

# --- 2. Insert the new 2022-Q4 row into "总计" (shifting the rest down) ----
$total.Rows.Item(2).Insert()

# The blank inserted row picks up stray formatting; reset it, then restore
# the row-index column's style (s="2", like every other data row) by
# copying the format from the row below.
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122) # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 2.32

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 5
$total.Range("D3").Value = 2.21

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 2.32

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 4
$total.Range("D5").Value = 0.38

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 6
$total.Range("D6").Value = 1.84

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 5
$total.Range("D7").Value = 2.29

$total.Range("A8").Value = 6
$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 5
$total.Range("D8").Value = 2.3
